$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Sheet2")

$rowCount = 93
$data = New-Object 'object[,]' $rowCount,3
$data[0,0] = 1
$data[0,1] = '黑暗之女'
$data[0,2] = 'Annie'
$data[1,0] = 2
$data[1,1] = '狂战士'
$data[1,2] = 'Olaf'
$data[2,0] = 3
$data[2,1] = '正义巨像'
$data[2,2] = 'Galio'
$data[3,0] = 4
$data[3,1] = '卡牌大师'
$data[3,2] = 'Twisted Fate'
$data[4,0] = 5
$data[4,1] = '德邦总管'
$data[4,2] = 'Xin Zhao'
$data[5,0] = 6
$data[5,1] = '无畏战车'
$data[5,2] = 'Urgot'
$data[6,0] = 8
$data[6,1] = '猩红收割者'
$data[6,2] = 'Vladimir'
$data[7,0] = 10
$data[7,1] = '正义天使'
$data[7,2] = 'Kayle'
$data[8,0] = 11
$data[8,1] = '无极剑圣'
$data[8,2] = 'Master Yi'
$data[9,0] = 12
$data[9,1] = '牛头酋长'
$data[9,2] = 'Alistar'
$data[10,0] = 13
$data[10,1] = '符文法师'
$data[10,2] = 'Ryze'
$data[11,0] = 15
$data[11,1] = '战争女神'
$data[11,2] = 'Sivir'
$data[12,0] = 16
$data[12,1] = '众星之子'
$data[12,2] = 'Soraka'
$data[13,0] = 17
$data[13,1] = '迅捷斥候'
$data[13,2] = 'Teemo'
$data[14,0] = 18
$data[14,1] = '麦林炮手'
$data[14,2] = 'Tristana'
$data[15,0] = 19
$data[15,1] = '祖安怒兽'
$data[15,2] = 'Warwick'
$data[16,0] = 21
$data[16,1] = '赏金猎人'
$data[16,2] = 'Miss Fortune'
$data[17,0] = 22
$data[17,1] = '寒冰射手'
$data[17,2] = 'Ashe'
$data[18,0] = 23
$data[18,1] = '蛮族之王'
$data[18,2] = 'Tryndamere'
$data[19,0] = 24
$data[19,1] = '武器大师'
$data[19,2] = 'Jax'
$data[20,0] = 25
$data[20,1] = '堕落天使'
$data[20,2] = 'Morgana'
$data[21,0] = 26
$data[21,1] = '时光守护者'
$data[21,2] = 'Zilean'
$data[22,0] = 29
$data[22,1] = '瘟疫之源'
$data[22,2] = 'Twitch'
$data[23,0] = 30
$data[23,1] = '死亡颂唱者'
$data[23,2] = 'Karthus'
$data[24,0] = 31
$data[24,1] = '虚空恐惧'
$data[24,2] = 'Cho''Gath'
$data[25,0] = 32
$data[25,1] = '殇之木乃伊'
$data[25,2] = 'Amumu'
$data[26,0] = 33
$data[26,1] = '披甲龙龟'
$data[26,2] = 'Rammus'
$data[27,0] = 36
$data[27,1] = '祖安狂人'
$data[27,2] = 'Dr. Mundo'
$data[28,0] = 37
$data[28,1] = '琴瑟仙女'
$data[28,2] = 'Sona'
$data[29,0] = 38
$data[29,1] = '虚空行者'
$data[29,2] = 'Kassadin'
$data[30,0] = 39
$data[30,1] = '刀锋舞者'
$data[30,2] = 'Irelia'
$data[31,0] = 40
$data[31,1] = '风暴之怒'
$data[31,2] = 'Janna'
$data[32,0] = 44
$data[32,1] = '瓦洛兰之盾'
$data[32,2] = 'Taric'
$data[33,0] = 45
$data[33,1] = '邪恶小法师'
$data[33,2] = 'Veigar'
$data[34,0] = 48
$data[34,1] = '巨魔之王'
$data[34,2] = 'Trundle'
$data[35,0] = 51
$data[35,1] = '皮城女警'
$data[35,2] = 'Caitlyn'
$data[36,0] = 53
$data[36,1] = '蒸汽机器人'
$data[36,2] = 'Blitzcrank'
$data[37,0] = 54
$data[37,1] = '熔岩巨兽'
$data[37,2] = 'Malphite'
$data[38,0] = 55
$data[38,1] = '不祥之刃'
$data[38,2] = 'Katarina'
$data[39,0] = 56
$data[39,1] = '永恒梦魇'
$data[39,2] = 'Nocturne'
$data[40,0] = 57
$data[40,1] = '扭曲树精'
$data[40,2] = 'Maokai'
$data[41,0] = 58
$data[41,1] = '荒漠屠夫'
$data[41,2] = 'Renekton'
$data[42,0] = 61
$data[42,1] = '发条魔灵'
$data[42,2] = 'Orianna'
$data[43,0] = 62
$data[43,1] = '齐天大圣'
$data[43,2] = 'Monkey King'
$data[44,0] = 63
$data[44,1] = '复仇焰魂'
$data[44,2] = 'Brand'
$data[45,0] = 67
$data[45,1] = '暗夜猎手'
$data[45,2] = 'Vayne'
$data[46,0] = 68
$data[46,1] = '机械公敌'
$data[46,2] = 'Rumble'
$data[47,0] = 69
$data[47,1] = '魔蛇之拥'
$data[47,2] = 'Cassiopeia'
$data[48,0] = 74
$data[48,1] = '大发明家'
$data[48,2] = 'Heimerdinger'
$data[49,0] = 75
$data[49,1] = '沙漠死神'
$data[49,2] = 'Nasus'
$data[50,0] = 76
$data[50,1] = '狂野女猎手'
$data[50,2] = 'Nidalee'
$data[51,0] = 78
$data[51,1] = '圣锤之毅'
$data[51,2] = 'Poppy'
$data[52,0] = 81
$data[52,1] = '探险家'
$data[52,2] = 'Ezreal'
$data[53,0] = 82
$data[53,1] = '铁铠冥魂'
$data[53,2] = 'Mordekaiser'
$data[54,0] = 83
$data[54,1] = '牧魂人'
$data[54,2] = 'Yorick'
$data[55,0] = 84
$data[55,1] = '离群之刺'
$data[55,2] = 'Akali'
$data[56,0] = 85
$data[56,1] = '狂暴之心'
$data[56,2] = 'Kennen'
$data[57,0] = 86
$data[57,1] = '德玛西亚之力'
$data[57,2] = 'Garen'
$data[58,0] = 89
$data[58,1] = '曙光女神'
$data[58,2] = 'Leona'
$data[59,0] = 90
$data[59,1] = '虚空先知'
$data[59,2] = 'Malzahar'
$data[60,0] = 96
$data[60,1] = '深渊巨口'
$data[60,2] = 'Kog''Maw'
$data[61,0] = 98
$data[61,1] = '暮光之眼'
$data[61,2] = 'Shen'
$data[62,0] = 99
$data[62,1] = '光辉女郎'
$data[62,2] = 'Lux'
$data[63,0] = 102
$data[63,1] = '龙血武姬'
$data[63,2] = 'Shyvana'
$data[64,0] = 103
$data[64,1] = '九尾妖狐'
$data[64,2] = 'Ahri'
$data[65,0] = 104
$data[65,1] = '法外狂徒'
$data[65,2] = 'Graves'
$data[66,0] = 105
$data[66,1] = '潮汐海灵'
$data[66,2] = 'Fizz'
$data[67,0] = 106
$data[67,1] = '不灭狂雷'
$data[67,2] = 'Volibear'
$data[68,0] = 110
$data[68,1] = '惩戒之箭'
$data[68,2] = 'Varus'
$data[69,0] = 111
$data[69,1] = '深海泰坦'
$data[69,2] = 'Nautilus'
$data[70,0] = 112
$data[70,1] = '机械先驱'
$data[70,2] = 'Viktor'
$data[71,0] = 113
$data[71,1] = '北地之怒'
$data[71,2] = 'Sejuani'
$data[72,0] = 115
$data[72,1] = '爆破鬼才'
$data[72,2] = 'Ziggs'
$data[73,0] = 117
$data[73,1] = '仙灵女巫'
$data[73,2] = 'Lulu'
$data[74,0] = 120
$data[74,1] = '战争之影'
$data[74,2] = 'Hecarim'
$data[75,0] = 121
$data[75,1] = '虚空掠夺者'
$data[75,2] = 'Khazix'
$data[76,0] = 122
$data[76,1] = '诺克萨斯之手'
$data[76,2] = 'Darius'
$data[77,0] = 127
$data[77,1] = '冰霜女巫'
$data[77,2] = 'Lissandra'
$data[78,0] = 131
$data[78,1] = '皎月女神'
$data[78,2] = 'Diana'
$data[79,0] = 143
$data[79,1] = '荆棘之兴'
$data[79,2] = 'Zyra'
$data[80,0] = 145
$data[80,1] = '虚空之女'
$data[80,2] = 'Kai''Sa'
$data[81,0] = 147
$data[81,1] = '星籁歌姬'
$data[81,2] = 'Seraphine'
$data[82,0] = 201
$data[82,1] = '弗雷尔卓德之心'
$data[82,2] = 'Braum'
$data[83,0] = 236
$data[83,1] = '圣枪游侠'
$data[83,2] = 'Lucian'
$data[84,0] = 254
$data[84,1] = '皮城执法官'
$data[84,2] = 'Vi'
$data[85,0] = 266
$data[85,1] = '暗裔剑魔'
$data[85,2] = 'Aatrox'
$data[86,0] = 497
$data[86,1] = '幻翎'
$data[86,2] = 'Rakan'
$data[87,0] = 711
$data[87,1] = '愁云使者'
$data[87,2] = 'Vex'
$data[88,0] = 875
$data[88,1] = '腕豪'
$data[88,2] = 'Sett'
$data[89,0] = 876
$data[89,1] = '含羞蓓蕾'
$data[89,2] = 'Lillia'
$data[90,0] = 888
$data[90,1] = '炼金男爵'
$data[90,2] = 'Renata Glasc'
$data[91,0] = 901
$data[91,1] = '炽炎雏龙'
$data[91,2] = 'Smolder'
$data[92,0] = 902
$data[92,1] = '明烛'
$data[92,2] = 'Milio'

$ws2.Range("A2:C94").Value2 = $data
